$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.071.26"
$ws.Range("E2").Value = "  -1.89%  "
$ws.Range("D3").Value = "2.495.87"
$ws.Range("E3").Value = "  -3.52%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'301.86"
$ws.Range("E5").Value = "  -0.38%  "
$ws.Range("D6").Value = "'94.95"
$ws.Range("E6").Value = "  -2.07%  "
$ws.Range("E7").Value = "  +0.73%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "'0.528"
$ws.Range("E9").Value = "  -4.13%  "
$ws.Range("D10").Value = "'36.13"
$ws.Range("E10").Value = "  -1.24%  "
$ws.Range("D11").Value = "'0.0803"
$ws.Range("E11").Value = "  -0.80%  "
$ws.Range("E12").Value = "  -2.10%  "
$ws.Range("E13").Value = "  -4.09%  "
$ws.Range("D14").Value = "2.875.35"
$ws.Range("E14").Value = "  -3.56%  "
$ws.Range("D15").Value = "2.492.02"
$ws.Range("E15").Value = "  -4.53%  "
$ws.Range("D16").Value = "'14.95"
$ws.Range("E16").Value = "  +3.88%  "
$ws.Range("D17").Value = "'0.845"
$ws.Range("E17").Value = "  -4.42%  "
$ws.Range("D18").Value = "42.086.32"
$ws.Range("E18").Value = "  -2.08%  "
$ws.Range("D19").Value = "'12.76"
$ws.Range("E19").Value = "  -1.15%  "
$ws.Range("E20").Value = "  -2.70%  "
$ws.Range("D21").Value = "'6.37"
$ws.Range("E21").Value = "  -4.46%  "
$ws.Range("D22").Value = "'70.73"
$ws.Range("E22").Value = "  -1.77%  "
$ws.Range("D23").Value = "'247.86"
$ws.Range("E23").Value = "  -2.76%  "
$ws.Range("D24").Value = "'2.88"
$ws.Range("E24").Value = "  -2.69%  "
$ws.Range("E25").Value = "  -6.18%  "
$ws.Range("D26").Value = "'26.42"
$ws.Range("E26").Value = "  -7.89%  "
$ws.Range("E27").Value = "  -0.15%  "
$ws.Range("E28").Value = "  +8.27%  "
$ws.Range("D29").Value = "'10.07"
$ws.Range("E29").Value = "  -1.38%  "
$ws.Range("D30").Value = "'37.08"
$ws.Range("E30").Value = "  -5.79%  "
$ws.Range("D31").Value = "'5.87"
$ws.Range("E31").Value = "  -2.86%  "
$ws.Range("D32").Value = "'154.30"
$ws.Range("E32").Value = "  -0.88%  "
$ws.Range("D33").Value = "'3.28"
$ws.Range("E33").Value = "  -2.77%  "
$ws.Range("E34").Value = "  -5.02%  "
$ws.Range("D35").Value = "'0.0778"
$ws.Range("E35").Value = "  -4.36%  "
$ws.Range("D36").Value = "'2.05"
$ws.Range("E36").Value = "  -6.29%  "
$ws.Range("D37").Value = "'18.41"
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("E38").Value = "  -0.92%  "
$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").Value = "'0.118"
$ws.Range("E39").Value = "  -1.46%  "
$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").Value = "'23.85"
$ws.Range("E40").Value = "  +2.27%  "
$ws.Range("D41").Value = "'3.82"
$ws.Range("E41").Value = "  -2.07%  "
$ws.Range("E42").Value = "  -2.27%  "
$ws.Range("D43").Value = "'0.998"
$ws.Range("D44").Value = "2.040.26"
$ws.Range("E44").Value = "  -1.47%  "
$ws.Range("E45").Value = "  -4.97%  "
$ws.Range("D46").Value = "'1.95"
$ws.Range("E46").Value = "  -6.70%  "
$ws.Range("D47").Value = "'8.91"
$ws.Range("E47").Value = "  -3.55%  "
$ws.Range("D48").Value = "'83.12"
$ws.Range("E48").Value = "  -2.47%  "
$ws.Range("D49").Value = "2.734.78"
$ws.Range("E49").Value = "  -3.55%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "'100.29"
$ws.Range("E50").Value = "  -5.67%  "
$ws.Range("B51").Value = "ordi"
$ws.Range("C51").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D51").Value = "'71.13"
$ws.Range("E51").Value = "  -7.14%  "
